$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Need to make music mute-able" bug entry (row 2) has been resolved
# (dev mode now mutes all sounds after escape is pressed once to exit a
# map), so remove that entry from the bug/request list.
$ws.Range("A2").ClearContents()

# Move the active selection to A2, as reflected in the saved view state.
$ws.Range("A2").Select()
